$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.943.75'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.380.40'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.54'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.10'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.387'
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.958.25'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.385.79'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.061.10'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.08'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.52'
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.17'
$ws.Range('E21').Value = '  -1.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.50'
$ws.Range('E22').Value = '  +3.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.548'
$ws.Range('E23').Value = '  -0.96%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.520.17'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.191'
$ws.Range('E27').Value = '  +7.42%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.20'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.94'
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('E33').Value = '  -3.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.18'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.91'
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.69'
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.415.26'
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.95'
$ws.Range('E38').Value = '  +0.73%  '
$ws.Range('E39').Value = '  -2.57%  '
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '25.96'
$ws.Range('E41').Value = '  -3.84%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.779'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.34'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('E45').Value = '  -2.64%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.429.98'
$ws.Range('E47').Value = '  -3.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.74'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('E49').Value = '  -2.03%  '
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('E51').Value = '  +7.65%  '
